$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Delete row 2 (value 10), shifting row 3 (value 13) up to become the new row 2
$ws.Rows.Item(2).Delete()

# Update the selection/active cell as recorded in the saved view state
$ws.Range("D12").Select()
